# Refresh crypto symbol list: update Price (D) and Volume(1h) (E) columns
# for the latest scrape snapshot. Values are written as literal text so
# percentages/prices keep their original formatting (e.g. "-4.56%"),
# matching how the source data was produced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextValue "D2" "307.42"
Set-TextValue "E2" "-4.56%"
Set-TextValue "D3" "39.82"
Set-TextValue "E3" "-7.50%"
Set-TextValue "D4" "5.120"
Set-TextValue "E4" "-1.51%"
Set-TextValue "D5" "0.07711"
Set-TextValue "E5" "-5.89%"
Set-TextValue "E6" "-1.88%"
Set-TextValue "D7" "1.629"
Set-TextValue "E7" "-11.35%"
Set-TextValue "D8" "0.8798"
Set-TextValue "E8" "-6.00%"
Set-TextValue "D9" "0.09986"
Set-TextValue "E9" "-10.31%"
Set-TextValue "D10" "0.1752"
Set-TextValue "E10" "-5.92%"
Set-TextValue "D11" "0.08952"
Set-TextValue "E11" "-5.46%"
Set-TextValue "D12" "0.04399"
Set-TextValue "E12" "-4.74%"
Set-TextValue "E13" "-0.29%"
Set-TextValue "D14" "0.001266"
Set-TextValue "E14" "-1.82%"
Set-TextValue "D15" "0.005841"
Set-TextValue "E15" "3.09%"
Set-TextValue "E16" "-0.20%"
Set-TextValue "D17" "2.438"
Set-TextValue "E17" "-3.26%"
Set-TextValue "E18" "-0.50%"
Set-TextValue "D19" "7.028"
Set-TextValue "E19" "-5.28%"
Set-TextValue "D21" "0.2850"
Set-TextValue "E21" "8.73%"
Set-TextValue "D22" "0.04156"
Set-TextValue "E22" "-0.13%"
Set-TextValue "E23" "-4.01%"
Set-TextValue "D24" "0.004054"
Set-TextValue "E24" "-6.06%"
Set-TextValue "D25" "0.0001219"
Set-TextValue "E25" "10.92%"
Set-TextValue "E26" "0.17%"
Set-TextValue "D38" "0.02346"
Set-TextValue "E38" "-13.73%"
Set-TextValue "D39" "0.05151"
Set-TextValue "E39" "-6.98%"
Set-TextValue "D40" "0.007912"
Set-TextValue "E40" "-0.73%"
Set-TextValue "D41" "0.1324"
Set-TextValue "E41" "-5.14%"
Set-TextValue "E42" "-3.36%"
Set-TextValue "D43" "0.001934"
Set-TextValue "E43" "-7.50%"
Set-TextValue "D44" "0.008513"
Set-TextValue "E44" "13.66%"
Set-TextValue "D45" "0.3064"
Set-TextValue "E45" "-4.45%"
Set-TextValue "D46" "0.00006519"
Set-TextValue "E46" "-6.55%"
Set-TextValue "E47" "0.07%"
Set-TextValue "E48" "98.64%"
Set-TextValue "D49" "0.006805"
Set-TextValue "E49" "96.55%"
Set-TextValue "E50" "0.07%"
Set-TextValue "E51" "0.07%"
